# Update Hbegf-Egfr LR-pairs worksheet with newly computed TPM-based values.
# Only the numeric metric columns (G-J, K-T as applicable) change; identifier
# columns (A-F) remain untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 8.459557
$ws.Range("H2").Value = 25.378671
$ws.Range("I2").Value = 0.3030961495696597
$ws.Range("J2").Value = 0.3030961495696597
$ws.Range("M2").Value = 1.701929666666667
$ws.Range("N2").Value = 5.105789
$ws.Range("O2").Value = 0.02105622887134972
$ws.Range("P2").Value = 0.02105622887134972
$ws.Range("Q2").Value = 14.39757102515767
$ws.Range("R2").Value = 129.578139226419
$ws.Range("S2").Value = 0.006382061895363601
$ws.Range("T2").Value = 0.006382061895363601
$ws.Range("G3").Value = 8.459557
$ws.Range("H3").Value = 25.378671
$ws.Range("I3").Value = 0.3030961495696597
$ws.Range("J3").Value = 0.3030961495696597
$ws.Range("O3").Value = 0.7732971809418951
$ws.Range("P3").Value = 0.7732971809418953
$ws.Range("Q3").Value = 528.7557023714791
$ws.Range("R3").Value = 4758.801321343311
$ws.Range("S3").Value = 0.2343833980165609
$ws.Range("T3").Value = 0.2343833980165609
$ws.Range("G4").Value = 8.459557
$ws.Range("H4").Value = 25.378671
$ws.Range("I4").Value = 0.3030961495696597
$ws.Range("J4").Value = 0.3030961495696597
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.4338690000000001
$ws.Range("N4").Value = 1.301607
$ws.Range("O4").Value = 0.005367815805265532
$ws.Range("P4").Value = 0.005367815805265533
$ws.Range("Q4").Value = 3.670339536033
$ws.Range("R4").Value = 33.03305582429701
$ws.Range("S4").Value = 0.001626964302175145
$ws.Range("T4").Value = 0.001626964302175145
$ws.Range("G5").Value = 8.459557
$ws.Range("H5").Value = 25.378671
$ws.Range("I5").Value = 0.3030961495696597
$ws.Range("J5").Value = 0.3030961495696597
$ws.Range("M5").Value = 15.972384
$ws.Range("N5").Value = 47.917152
$ws.Range("O5").Value = 0.1976099128607259
$ws.Range("P5").Value = 0.1976099128607259
$ws.Range("Q5").Value = 135.119292873888
$ws.Range("R5").Value = 1216.073635864992
$ws.Range("S5").Value = 0.05989480370488201
$ws.Range("T5").Value = 0.059894803704882
$ws.Range("G6").Value = 8.459557
$ws.Range("H6").Value = 25.378671
$ws.Range("I6").Value = 0.3030961495696597
$ws.Range("J6").Value = 0.3030961495696597
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.2157183333333333
$ws.Range("N6").Value = 0.647155
$ws.Range("O6").Value = 0.002668861520763652
$ws.Range("P6").Value = 0.002668861520763652
$ws.Range("Q6").Value = 1.824881536778334
$ws.Range("R6").Value = 16.423933831005
$ws.Range("S6").Value = 0.0008089216506780896
$ws.Range("T6").Value = 0.0008089216506780894
$ws.Range("I7").Value = 0.4601547065605718
$ws.Range("J7").Value = 0.4601547065605718
$ws.Range("M7").Value = 1.701929666666667
$ws.Range("N7").Value = 5.105789
$ws.Range("O7").Value = 0.02105622887134972
$ws.Range("P7").Value = 0.02105622887134972
$ws.Range("Q7").Value = 21.85811360412478
$ws.Range("R7").Value = 196.723022437123
$ws.Range("S7").Value = 0.00968912281756817
$ws.Range("T7").Value = 0.009689122817568169
$ws.Range("I8").Value = 0.4601547065605718
$ws.Range("J8").Value = 0.4601547065605718
$ws.Range("O8").Value = 0.7732971809418951
$ws.Range("P8").Value = 0.7732971809418953
$ws.Range("S8").Value = 0.3558363373804352
$ws.Range("T8").Value = 0.3558363373804352
$ws.Range("I9").Value = 0.4601547065605718
$ws.Range("J9").Value = 0.4601547065605718
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.4338690000000001
$ws.Range("N9").Value = 1.301607
$ws.Range("O9").Value = 0.005367815805265532
$ws.Range("P9").Value = 0.005367815805265533
$ws.Range("Q9").Value = 5.572238428561
$ws.Range("R9").Value = 50.150145857049
$ws.Range("S9").Value = 0.00247002570674316
$ws.Range("T9").Value = 0.002470025706743161
$ws.Range("I10").Value = 0.4601547065605718
$ws.Range("J10").Value = 0.4601547065605718
$ws.Range("M10").Value = 15.972384
$ws.Range("N10").Value = 47.917152
$ws.Range("O10").Value = 0.1976099128607259
$ws.Range("P10").Value = 0.1976099128607259
$ws.Range("Q10").Value = 205.135494632096
$ws.Range("R10").Value = 1846.219451688864
$ws.Range("S10").Value = 0.0909311314658875
$ws.Range("T10").Value = 0.09093113146588751
$ws.Range("I11").Value = 0.4601547065605718
$ws.Range("J11").Value = 0.4601547065605718
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 0.2157183333333333
$ws.Range("N11").Value = 0.647155
$ws.Range("O11").Value = 0.002668861520763652
$ws.Range("P11").Value = 0.002668861520763652
$ws.Range("Q11").Value = 2.770499820787222
$ws.Range("R11").Value = 24.934498387085
$ws.Range("S11").Value = 0.0012280891899378
$ws.Range("T11").Value = 0.0012280891899378
$ws.Range("G12").Value = 1.955432333333333
$ws.Range("H12").Value = 5.866296999999999
$ws.Range("I12").Value = 0.0700608803720276
$ws.Range("J12").Value = 0.0700608803720276
$ws.Range("M12").Value = 1.701929666666667
$ws.Range("N12").Value = 5.105789
$ws.Range("O12").Value = 0.02105622887134972
$ws.Range("P12").Value = 0.02105622887134972
$ws.Range("Q12").Value = 3.328008299259222
$ws.Range("R12").Value = 29.952074693333
$ws.Range("S12").Value = 0.001475217932041666
$ws.Range("T12").Value = 0.001475217932041666
$ws.Range("G13").Value = 1.955432333333333
$ws.Range("H13").Value = 5.866296999999999
$ws.Range("I13").Value = 0.0700608803720276
$ws.Range("J13").Value = 0.0700608803720276
$ws.Range("O13").Value = 0.7732971809418951
$ws.Range("P13").Value = 0.7732971809418953
$ws.Range("Q13").Value = 122.222238924753
$ws.Range("R13").Value = 1100.000150322777
$ws.Range("S13").Value = 0.05417788128599629
$ws.Range("T13").Value = 0.0541778812859963
$ws.Range("G14").Value = 1.955432333333333
$ws.Range("H14").Value = 5.866296999999999
$ws.Range("I14").Value = 0.0700608803720276
$ws.Range("J14").Value = 0.0700608803720276
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 0.4338690000000001
$ws.Range("N14").Value = 1.301607
$ws.Range("O14").Value = 0.005367815805265532
$ws.Range("P14").Value = 0.005367815805265533
$ws.Range("Q14").Value = 0.848401471031
$ws.Range("R14").Value = 7.635613239279
$ws.Range("S14").Value = 0.0003760739009917874
$ws.Range("T14").Value = 0.0003760739009917875
$ws.Range("G15").Value = 1.955432333333333
$ws.Range("H15").Value = 5.866296999999999
$ws.Range("I15").Value = 0.0700608803720276
$ws.Range("J15").Value = 0.0700608803720276
$ws.Range("M15").Value = 15.972384
$ws.Range("N15").Value = 47.917152
$ws.Range("O15").Value = 0.1976099128607259
$ws.Range("P15").Value = 0.1976099128607259
$ws.Range("Q15").Value = 31.232916114016
$ws.Range("R15").Value = 281.096245026144
$ws.Range("S15").Value = 0.01384472446526211
$ws.Range("T15").Value = 0.01384472446526212
$ws.Range("G16").Value = 1.955432333333333
$ws.Range("H16").Value = 5.866296999999999
$ws.Range("I16").Value = 0.0700608803720276
$ws.Range("J16").Value = 0.0700608803720276
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 0.2157183333333333
$ws.Range("N16").Value = 0.647155
$ws.Range("O16").Value = 0.002668861520763652
$ws.Range("P16").Value = 0.002668861520763652
$ws.Range("Q16").Value = 0.4218226038927778
$ws.Range("R16").Value = 3.796403435035
$ws.Range("S16").Value = 0.0001869827877357299
$ws.Range("T16").Value = 0.0001869827877357299
$ws.Range("G17").Value = 2.929608
$ws.Range("H17").Value = 8.788824000000002
$ws.Range("I17").Value = 0.1049644685352285
$ws.Range("J17").Value = 0.1049644685352285
$ws.Range("M17").Value = 1.701929666666667
$ws.Range("N17").Value = 5.105789
$ws.Range("O17").Value = 0.02105622887134972
$ws.Range("P17").Value = 0.02105622887134972
$ws.Range("Q17").Value = 4.985986766904
$ws.Range("R17").Value = 44.87388090213601
$ws.Range("S17").Value = 0.002210155872837357
$ws.Range("T17").Value = 0.002210155872837357
$ws.Range("G18").Value = 2.929608
$ws.Range("H18").Value = 8.788824000000002
$ws.Range("I18").Value = 0.1049644685352285
$ws.Range("J18").Value = 0.1049644685352285
$ws.Range("O18").Value = 0.7732971809418951
$ws.Range("P18").Value = 0.7732971809418953
$ws.Range("Q18").Value = 183.112063162776
$ws.Range("R18").Value = 1648.008568464984
$ws.Range("S18").Value = 0.08116872761735644
$ws.Range("T18").Value = 0.08116872761735645
$ws.Range("G19").Value = 2.929608
$ws.Range("H19").Value = 8.788824000000002
$ws.Range("I19").Value = 0.1049644685352285
$ws.Range("J19").Value = 0.1049644685352285
$ws.Range("K19").Value = 3
$ws.Range("L19").Value = 1
$ws.Range("M19").Value = 0.4338690000000001
$ws.Range("N19").Value = 1.301607
$ws.Range("O19").Value = 0.005367815805265532
$ws.Range("P19").Value = 0.005367815805265533
$ws.Range("Q19").Value = 1.271066093352
$ws.Range("R19").Value = 11.439594840168
$ws.Range("S19").Value = 0.000563429933194696
$ws.Range("T19").Value = 0.0005634299331946962
$ws.Range("G20").Value = 2.929608
$ws.Range("H20").Value = 8.788824000000002
$ws.Range("I20").Value = 0.1049644685352285
$ws.Range("J20").Value = 0.1049644685352285
$ws.Range("M20").Value = 15.972384
$ws.Range("N20").Value = 47.917152
$ws.Range("O20").Value = 0.1976099128607259
$ws.Range("P20").Value = 0.1976099128607259
$ws.Range("Q20").Value = 46.79282394547201
$ws.Range("R20").Value = 421.1354155092481
$ws.Range("S20").Value = 0.02074201948071891
$ws.Range("T20").Value = 0.02074201948071891
$ws.Range("G21").Value = 2.929608
$ws.Range("H21").Value = 8.788824000000002
$ws.Range("I21").Value = 0.1049644685352285
$ws.Range("J21").Value = 0.1049644685352285
$ws.Range("K21").Value = 3
$ws.Range("L21").Value = 1
$ws.Range("M21").Value = 0.2157183333333333
$ws.Range("N21").Value = 0.647155
$ws.Range("O21").Value = 0.002668861520763652
$ws.Range("P21").Value = 0.002668861520763652
$ws.Range("Q21").Value = 0.6319701550800001
$ws.Range("R21").Value = 5.687731395720001
$ws.Range("S21").Value = 0.0002801356311210785
$ws.Range("T21").Value = 0.0002801356311210785
$ws.Range("G22").Value = 1.722740333333333
$ws.Range("H22").Value = 5.168221
$ws.Range("I22").Value = 0.06172379496251228
$ws.Range("J22").Value = 0.06172379496251227
$ws.Range("M22").Value = 1.701929666666667
$ws.Range("N22").Value = 5.105789
$ws.Range("O22").Value = 0.02105622887134972
$ws.Range("P22").Value = 0.02105622887134972
$ws.Range("Q22").Value = 2.931982881263222
$ws.Range("R22").Value = 26.387845931369
$ws.Range("S22").Value = 0.001299670353538921
$ws.Range("T22").Value = 0.001299670353538921
$ws.Range("G23").Value = 1.722740333333333
$ws.Range("H23").Value = 5.168221
$ws.Range("I23").Value = 0.06172379496251228
$ws.Range("J23").Value = 0.06172379496251227
$ws.Range("O23").Value = 0.7732971809418951
$ws.Range("P23").Value = 0.7732971809418953
$ws.Range("Q23").Value = 107.678070489429
$ws.Range("R23").Value = 969.102634404861
$ws.Range("S23").Value = 0.04773083664154629
$ws.Range("T23").Value = 0.04773083664154629
$ws.Range("G24").Value = 1.722740333333333
$ws.Range("H24").Value = 5.168221
$ws.Range("I24").Value = 0.06172379496251228
$ws.Range("J24").Value = 0.06172379496251227
$ws.Range("K24").Value = 3
$ws.Range("L24").Value = 1
$ws.Range("M24").Value = 0.4338690000000001
$ws.Range("N24").Value = 1.301607
$ws.Range("O24").Value = 0.005367815805265532
$ws.Range("P24").Value = 0.005367815805265533
$ws.Range("Q24").Value = 0.747443625683
$ws.Range("R24").Value = 6.726992631147001
$ws.Range("S24").Value = 0.0003313219621607424
$ws.Range("T24").Value = 0.0003313219621607424
$ws.Range("G25").Value = 1.722740333333333
$ws.Range("H25").Value = 5.168221
$ws.Range("I25").Value = 0.06172379496251228
$ws.Range("J25").Value = 0.06172379496251227
$ws.Range("M25").Value = 15.972384
$ws.Range("N25").Value = 47.917152
$ws.Range("O25").Value = 0.1976099128607259
$ws.Range("P25").Value = 0.1976099128607259
$ws.Range("Q25").Value = 27.516270136288
$ws.Range("R25").Value = 247.646431226592
$ws.Range("S25").Value = 0.01219723374397536
$ws.Range("T25").Value = 0.01219723374397536
$ws.Range("G26").Value = 1.722740333333333
$ws.Range("H26").Value = 5.168221
$ws.Range("I26").Value = 0.06172379496251228
$ws.Range("J26").Value = 0.06172379496251227
$ws.Range("K26").Value = 3
$ws.Range("L26").Value = 1
$ws.Range("M26").Value = 0.2157183333333333
$ws.Range("N26").Value = 0.647155
$ws.Range("O26").Value = 0.002668861520763652
$ws.Range("P26").Value = 0.002668861520763652
$ws.Range("Q26").Value = 0.3716266734727778
$ws.Range("R26").Value = 3.344640061255
$ws.Range("S26").Value = 0.0001647322612909544
$ws.Range("T26").Value = 0.0001647322612909544
